$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 45.98144433333334
$ws.Range("H2").Value = 137.944333
$ws.Range("I2").Value = 0.9841234286873372
$ws.Range("J2").Value = 0.984123428687337
$ws.Range("M2").Value = 70.46291600000001
$ws.Range("N2").Value = 211.388748
$ws.Range("O2").Value = 0.5276750397950939
$ws.Range("P2").Value = 0.5276750397950939
$ws.Range("Q2").Value = 3239.986649618343
$ws.Range("R2").Value = 29159.87984656509
$ws.Range("S2").Value = 0.5192973693958749
$ws.Range("T2").Value = 0.5192973693958748
$ws.Range("G3").Value = 45.98144433333334
$ws.Range("H3").Value = 137.944333
$ws.Range("I3").Value = 0.9841234286873372
$ws.Range("J3").Value = 0.984123428687337
$ws.Range("O3").Value = 0.07361176802536967
$ws.Range("P3").Value = 0.07361176802536967
$ws.Range("Q3").Value = 451.9848915908826
$ws.Range("R3").Value = 4067.864024317943
$ws.Range("S3").Value = 0.0724430655408637
$ws.Range("T3").Value = 0.07244306554086369
$ws.Range("G4").Value = 45.98144433333334
$ws.Range("H4").Value = 137.944333
$ws.Range("I4").Value = 0.9841234286873372
$ws.Range("J4").Value = 0.984123428687337
$ws.Range("M4").Value = 42.505498
$ws.Range("N4").Value = 127.516494
$ws.Range("O4").Value = 0.3183105613832428
$ws.Range("P4").Value = 0.3183105613832428
$ws.Range("Q4").Value = 1954.464190147611
$ws.Range("R4").Value = 17590.1777113285
$ws.Range("S4").Value = 0.313256881055868
$ws.Range("T4").Value = 0.313256881055868
$ws.Range("G5").Value = 45.98144433333334
$ws.Range("H5").Value = 137.944333
$ws.Range("I5").Value = 0.9841234286873372
$ws.Range("J5").Value = 0.984123428687337
$ws.Range("M5").Value = 10.73653933333333
$ws.Range("N5").Value = 32.209618
$ws.Range("O5").Value = 0.08040263079629371
$ws.Range("P5").Value = 0.08040263079629371
$ws.Range("Q5").Value = 493.6815856883105
$ws.Range("R5").Value = 4443.134271194794
$ws.Range("S5").Value = 0.07912611269473066
$ws.Range("T5").Value = 0.07912611269473065
$ws.Range("I6").Value = 0.002244435796517234
$ws.Range("J6").Value = 0.002244435796517234
$ws.Range("M6").Value = 70.46291600000001
$ws.Range("N6").Value = 211.388748
$ws.Range("O6").Value = 0.5276750397950939
$ws.Range("P6").Value = 0.5276750397950939
$ws.Range("Q6").Value = 7.389258099810667
$ws.Range("R6").Value = 66.50332289829601
$ws.Range("S6").Value = 0.001184332748244765
$ws.Range("T6").Value = 0.001184332748244765
$ws.Range("I7").Value = 0.002244435796517234
$ws.Range("J7").Value = 0.002244435796517234
$ws.Range("O7").Value = 0.07361176802536967
$ws.Range("P7").Value = 0.07361176802536967
$ws.Range("S7").Value = 0.0001652168872010625
$ws.Range("T7").Value = 0.0001652168872010624
$ws.Range("I8").Value = 0.002244435796517234
$ws.Range("J8").Value = 0.002244435796517234
$ws.Range("M8").Value = 42.505498
$ws.Range("N8").Value = 127.516494
$ws.Range("O8").Value = 0.3183105613832428
$ws.Range("P8").Value = 0.3183105613832428
$ws.Range("Q8").Value = 4.457438227265333
$ws.Range("R8").Value = 40.116944045388
$ws.Range("S8").Value = 0.0007144276183780467
$ws.Range("T8").Value = 0.0007144276183780465
$ws.Range("I9").Value = 0.002244435796517234
$ws.Range("J9").Value = 0.002244435796517234
$ws.Range("M9").Value = 10.73653933333333
$ws.Range("N9").Value = 32.209618
$ws.Range("O9").Value = 0.08040263079629371
$ws.Range("P9").Value = 0.08040263079629371
$ws.Range("Q9").Value = 1.125912249115111
$ws.Range("R9").Value = 10.133210242036
$ws.Range("S9").Value = 0.0001804585426933606
$ws.Range("T9").Value = 0.0001804585426933606
$ws.Range("G10").Value = 0.547937
$ws.Range("H10").Value = 1.643811
$ws.Range("I10").Value = 0.01172728797372169
$ws.Range("J10").Value = 0.01172728797372169
$ws.Range("M10").Value = 70.46291600000001
$ws.Range("N10").Value = 211.388748
$ws.Range("O10").Value = 0.5276750397950939
$ws.Range("P10").Value = 0.5276750397950939
$ws.Range("Q10").Value = 38.609238804292
$ws.Range("R10").Value = 347.483149238628
$ws.Range("S10").Value = 0.006188197148222121
$ws.Range("T10").Value = 0.006188197148222119
$ws.Range("G11").Value = 0.547937
$ws.Range("H11").Value = 1.643811
$ws.Range("I11").Value = 0.01172728797372169
$ws.Range("J11").Value = 0.01172728797372169
$ws.Range("O11").Value = 0.07361176802536967
$ws.Range("P11").Value = 0.07361176802536967
$ws.Range("Q11").Value = 5.386069296742333
$ws.Range("R11").Value = 48.474623670681
$ws.Range("S11").Value = 0.0008632664018883088
$ws.Range("T11").Value = 0.0008632664018883086
$ws.Range("G12").Value = 0.547937
$ws.Range("H12").Value = 1.643811
$ws.Range("I12").Value = 0.01172728797372169
$ws.Range("J12").Value = 0.01172728797372169
$ws.Range("M12").Value = 42.505498
$ws.Range("N12").Value = 127.516494
$ws.Range("O12").Value = 0.3183105613832428
$ws.Range("P12").Value = 0.3183105613832428
$ws.Range("Q12").Value = 23.290335057626
$ws.Range("R12").Value = 209.613015518634
$ws.Range("S12").Value = 0.003732919618418304
$ws.Range("T12").Value = 0.003732919618418303
$ws.Range("G13").Value = 0.547937
$ws.Range("H13").Value = 1.643811
$ws.Range("I13").Value = 0.01172728797372169
$ws.Range("J13").Value = 0.01172728797372169
$ws.Range("M13").Value = 10.73653933333333
$ws.Range("N13").Value = 32.209618
$ws.Range("O13").Value = 0.08040263079629371
$ws.Range("P13").Value = 0.08040263079629371
$ws.Range("Q13").Value = 5.882947152688667
$ws.Range("R13").Value = 52.946524374198
$ws.Range("S13").Value = 0.0009429048051929608
$ws.Range("T13").Value = 0.0009429048051929604
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.08900066666666666
$ws.Range("H14").Value = 0.267002
$ws.Range("I14").Value = 0.001904847542424061
$ws.Range("J14").Value = 0.001904847542424061
$ws.Range("M14").Value = 70.46291600000001
$ws.Range("N14").Value = 211.388748
$ws.Range("O14").Value = 0.5276750397950939
$ws.Range("P14").Value = 0.5276750397950939
$ws.Range("Q14").Value = 6.271246499277334
$ws.Range("R14").Value = 56.441218493496
$ws.Range("S14").Value = 0.001005140502752203
$ws.Range("T14").Value = 0.001005140502752203
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.08900066666666666
$ws.Range("H15").Value = 0.267002
$ws.Range("I15").Value = 0.001904847542424061
$ws.Range("J15").Value = 0.001904847542424061
$ws.Range("O15").Value = 0.07361176802536967
$ws.Range("P15").Value = 0.07361176802536967
$ws.Range("Q15").Value = 0.8748519594824443
$ws.Range("R15").Value = 7.873667635341999
$ws.Range("S15").Value = 0.0001402191954166155
$ws.Range("T15").Value = 0.0001402191954166155
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.08900066666666666
$ws.Range("H16").Value = 0.267002
$ws.Range("I16").Value = 0.001904847542424061
$ws.Range("J16").Value = 0.001904847542424061
$ws.Range("M16").Value = 42.505498
$ws.Range("N16").Value = 127.516494
$ws.Range("O16").Value = 0.3183105613832428
$ws.Range("P16").Value = 0.3183105613832428
$ws.Range("Q16").Value = 3.783017658998666
$ws.Range("R16").Value = 34.04715893098799
$ws.Range("S16").Value = 0.0006063330905784934
$ws.Range("T16").Value = 0.0006063330905784934
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.08900066666666666
$ws.Range("H17").Value = 0.267002
$ws.Range("I17").Value = 0.001904847542424061
$ws.Range("J17").Value = 0.001904847542424061
$ws.Range("M17").Value = 10.73653933333333
$ws.Range("N17").Value = 32.209618
$ws.Range("O17").Value = 0.08040263079629371
$ws.Range("P17").Value = 0.08040263079629371
$ws.Range("Q17").Value = 0.9555591583595554
$ws.Range("R17").Value = 8.600032425235998
$ws.Range("S17").Value = 0.0001531547536767493
$ws.Range("T17").Value = 0.0001531547536767492
